# Fruta / hortaliza, semanal
# A new daily price record (Región de O'Higgins, fecha 2023-02-14 / serial 44971)
# was inserted as row 20, pushing all subsequent records down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 20; every row from the old 20 onward
# (including their formatting) shifts down to 21.., exactly matching the
# diff where old row N's values now live in row N+1.
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new record's data.
$ws.Range("A20").Value = 9
$ws.Range("B20").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C20").Value = "Metropolitana"
$ws.Range("D20").Value = 44971
$ws.Range("E20").Value = 13
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100101
$ws.Range("H20").Value = "Berries"
$ws.Range("I20").Value = 100101004
$ws.Range("J20").Value = "Frambuesa"
$ws.Range("K20").Value = "Sin especificar"
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 400
$ws.Range("N20").Value = 5600
$ws.Range("O20").Value = 5600
$ws.Range("P20").Value = 5600
$ws.Range("Q20").Value = "`$/bandeja 2 kilos"
$ws.Range("R20").Value = "Región de O'Higgins"
$ws.Range("S20").Value = 2800
$ws.Range("T20").Value = 2
